$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu results for the 380 kV case (rows 2-25)
# Columns B:F and I:N are updated per row; column G (slack, always 1) and column H (blank) are left untouched.

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.030514092294297
$bf[0,2] = 1.038837503113702
$bf[0,3] = 1.030182920204198
$bf[0,4] = 1.04751285351414
$ws.Range("B2:F2").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032821483893285
$in[0,1] = 1.035654917450162
$in[0,2] = 1.04162439917533
$in[0,3] = 1.032994671139967
$in[0,4] = 1.050275282250775
$in[0,5] = 1.015786072005686
$ws.Range("I2:N2").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.031744729451569
$bf[0,2] = 1.039811002195391
$bf[0,3] = 1.031236808988346
$bf[0,4] = 1.048723766484677
$ws.Range("B3:F3").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033054304476766
$in[0,1] = 1.036525606524905
$in[0,2] = 1.04240756229441
$in[0,3] = 1.033856201305355
$in[0,4] = 1.051297020078623
$in[0,5] = 1.016082382429905
$ws.Range("I3:N3").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032540738531023
$bf[0,2] = 1.040440436168522
$bf[0,3] = 1.03191880723601
$bf[0,4] = 1.049507271287637
$ws.Range("B4:F4").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033203362106256
$in[0,1] = 1.037088229512234
$in[0,2] = 1.042913223081675
$in[0,3] = 1.034413136854964
$in[0,4] = 1.051957550637025
$in[0,5] = 1.016273654327489
$ws.Range("I4:N4").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032875312319155
$bf[0,2] = 1.040704935065175
$bf[0,3] = 1.032205535255764
$bf[0,4] = 1.049836649430074
$ws.Range("B5:F5").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033265644873019
$in[0,1] = 1.037324572869677
$in[0,2] = 1.043125540954408
$in[0,3] = 1.034647146294582
$in[0,4] = 1.052235094848545
$in[0,5] = 1.016353955038827
$ws.Range("I5:N5").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032931484784578
$bf[0,2] = 1.04074933886972
$bf[0,3] = 1.0322536791218
$bf[0,4] = 1.049891953087741
$ws.Range("B6:F6").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033276080085918
$in[0,1] = 1.037364245244046
$in[0,2] = 1.043161174728349
$in[0,3] = 1.034686430137258
$in[0,4] = 1.052281687405223
$in[0,5] = 1.016367431440526
$ws.Range("I6:N6").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032545209388359
$bf[0,2] = 1.040443970866968
$bf[0,3] = 1.031922638445681
$bf[0,4] = 1.049511672481163
$ws.Range("B7:F7").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033204195828606
$in[0,1] = 1.037091388262621
$in[0,2] = 1.04291606111146
$in[0,3] = 1.034416264195304
$in[0,4] = 1.051961259755973
$in[0,5] = 1.016274727741339
$ws.Range("I7:N7").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.030930054445627
$bf[0,2] = 1.039166602701117
$bf[0,3] = 1.030539074781126
$bf[0,4] = 1.047922095787203
$ws.Range("B8:F8").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032900496675647
$in[0,1] = 1.035949331306338
$in[0,2] = 1.041889300776343
$in[0,3] = 1.033285940135241
$in[0,4] = 1.050620708804761
$in[0,5] = 1.015886307093325
$ws.Range("I8:N8").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.028081578872156
$bf[0,2] = 1.036911957911112
$bf[0,3] = 1.028101478373122
$bf[0,4] = 1.045120706425106
$ws.Range("B9:F9").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032353128080892
$in[0,1] = 1.033930919244611
$in[0,2] = 1.040071563631953
$in[0,3] = 1.031290039273731
$in[0,4] = 1.048253815553779
$in[0,5] = 1.01519831715958
$ws.Range("I9:N9").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.026180857481521
$bf[0,2] = 1.035406269543669
$bf[0,3] = 1.026476619080684
$bf[0,4] = 1.043252776306461
$ws.Range("B10:F10").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031979986288949
$in[0,1] = 1.032581221750401
$in[0,2] = 1.038853990789093
$in[0,3] = 1.029956592247961
$in[0,4] = 1.046672660144735
$in[0,5] = 1.014737254048741
$ws.Range("I10:N10").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.025357378627392
$bf[0,2] = 1.034753659838803
$bf[0,3] = 1.0257730677122
$bf[0,4] = 1.042443838725482
$ws.Range("B11:F11").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031816455254422
$in[0,1] = 1.031995799105675
$in[0,2] = 1.038325389942416
$in[0,3] = 1.02937850292245
$in[0,4] = 1.045987217136144
$in[0,5] = 1.014537033646018
$ws.Range("I11:N11").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.025051430747514
$bf[0,2] = 1.03451115445568
$bf[0,3] = 1.025511739283771
$bf[0,4] = 1.042143344428073
$ws.Range("B12:F12").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031755417889844
$in[0,1] = 1.031778195912641
$in[0,2] = 1.038128834763067
$in[0,3] = 1.029163668399871
$in[0,4] = 1.045732492226368
$in[0,5] = 1.014462575709426
$ws.Range("I12:N12").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.025117060874557
$bf[0,2] = 1.034563177099646
$bf[0,3] = 1.025567795060392
$bf[0,4] = 1.042207802362406
$ws.Range("B13:F13").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031768523953066
$in[0,1] = 1.031824879403935
$in[0,2] = 1.03817100600947
$in[0,3] = 1.029209755974316
$in[0,4] = 1.045787137082634
$in[0,5] = 1.014478551141548
$ws.Range("I13:N13").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.0253320903492
$bf[0,2] = 1.034733616256877
$bf[0,3] = 1.025751466194064
$bf[0,4] = 1.042419000170121
$ws.Range("B14:F14").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031811415900431
$in[0,1] = 1.031977815051332
$in[0,2] = 1.038309146907489
$in[0,3] = 1.029360746813449
$in[0,4] = 1.045966163971675
$in[0,5] = 1.014530880710613
$ws.Range("I14:N14").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.025464567586029
$bf[0,2] = 1.034838616539907
$bf[0,3] = 1.025864632240624
$bf[0,4] = 1.042549123547263
$ws.Range("B15:F15").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.03183780398598
$in[0,1] = 1.032072023670154
$in[0,2] = 1.038394232299374
$in[0,3] = 1.029453763106337
$in[0,4] = 1.046076452277168
$in[0,5] = 1.014563111120946
$ws.Range("I15:N15").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.026235499243739
$bf[0,2] = 1.035449567564505
$bf[0,3] = 1.026523311821351
$bf[0,4] = 1.04330646030033
$ws.Range("B16:F16").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031990798009129
$in[0,1] = 1.032620053219746
$in[0,2] = 1.038889043013185
$in[0,3] = 1.02999494329547
$in[0,4] = 1.046718133891647
$in[0,5] = 1.014750529819864
$ws.Range("I16:N16").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.026718960677008
$bf[0,2] = 1.035832629534184
$bf[0,3] = 1.02693648923366
$bf[0,4] = 1.043781486295416
$ws.Range("B17:F17").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032086242531993
$in[0,1] = 1.032963550155783
$in[0,2] = 1.039199053244864
$in[0,3] = 1.030334223375566
$in[0,4] = 1.04712043041208
$in[0,5] = 1.014867937797974
$ws.Range("I17:N17").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.027000911954
$bf[0,2] = 1.036056001765119
$bf[0,3] = 1.027177491000912
$bf[0,4] = 1.044058550306281
$ws.Range("B18:F18").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032141724803698
$in[0,1] = 1.033163810247657
$in[0,2] = 1.039379743549689
$in[0,3] = 1.030532052468538
$in[0,4] = 1.047355006847257
$in[0,5] = 1.014936364242309
$ws.Range("I18:N18").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.027097042747464
$bf[0,2] = 1.036132155562385
$bf[0,3] = 1.027259666827179
$bf[0,4] = 1.044153020236618
$ws.Range("B19:F19").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032160610797
$in[0,1] = 1.033232077555635
$in[0,2] = 1.039441331685745
$in[0,3] = 1.030599495718503
$in[0,4] = 1.047434978459292
$in[0,5] = 1.014959686467583
$ws.Range("I19:N19").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.02666709437286
$bf[0,2] = 1.035791536965514
$bf[0,3] = 1.026892158975002
$bf[0,4] = 1.043730521616129
$ws.Range("B20:F20").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032076021780878
$in[0,1] = 1.032926706113841
$in[0,2] = 1.039165805891919
$in[0,3] = 1.030297828800637
$in[0,4] = 1.047077275710983
$in[0,5] = 1.014855346792193
$ws.Range("I20:N20").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.025268771516847
$bf[0,2] = 1.03468342886455
$bf[0,3] = 1.025697379584554
$bf[0,4] = 1.042356808217372
$ws.Range("B21:F21").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031798793428699
$in[0,1] = 1.031932783501441
$in[0,2] = 1.038268473644717
$in[0,3] = 1.0293162867245
$in[0,4] = 1.045913448365928
$in[0,5] = 1.014515473366766
$ws.Range("I21:N21").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.024389177206126
$bf[0,2] = 1.033986154927301
$bf[0,3] = 1.024946184391882
$bf[0,4] = 1.041492989901769
$ws.Range("B22:F22").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031622783933831
$in[0,1] = 1.031306989507024
$in[0,2] = 1.037703073669603
$in[0,3] = 1.02869853585115
$in[0,4] = 1.045181003770623
$in[0,5] = 1.014301276913624
$ws.Range("I22:N22").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.024855506955768
$bf[0,2] = 1.034355846820519
$bf[0,3] = 1.025344406777216
$bf[0,4] = 1.041950927485146
$ws.Range("B23:F23").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.031716251669606
$in[0,1] = 1.0316388182951
$in[0,2] = 1.038002918204844
$in[0,3] = 1.029026076224647
$in[0,4] = 1.045569353487125
$in[0,5] = 1.014414874489196
$ws.Range("I23:N23").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.026690530657818
$bf[0,2] = 1.035810105118088
$bf[0,3] = 1.026912189903455
$bf[0,4] = 1.043753550393805
$ws.Range("B24:F24").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032080640682077
$in[0,1] = 1.032943354646681
$in[0,2] = 1.039180829351891
$in[0,3] = 1.030314274150258
$in[0,4] = 1.0470967756987
$in[0,5] = 1.014861036297688
$ws.Range("I24:N24").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.028818273580461
$bf[0,2] = 1.037495289698752
$bf[0,3] = 1.028731613666942
$bf[0,4] = 1.045844984769203
$ws.Range("B25:F25").Value = $bf
$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.032496085080712
$in[0,1] = 1.034453441970432
$in[0,2] = 1.040542500682636
$in[0,3] = 1.031806523839726
$in[0,4] = 1.048866276574539
$in[0,5] = 1.015376601073277
$ws.Range("I25:N25").Value = $in

Write-Output "Updated vm_pu values for rows 2-25 (380 kV case)"